# Auto-generated edit script: updates Price (D) and Volume(1h) (E) columns
# on the crypto symbol list sheet, matching the upstream GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

# Row 2
Set-TextValue "D2" "305.77"
Set-TextValue "E2" "5.88%"

# Row 3
Set-TextValue "D3" "32.20"
Set-TextValue "E3" "9.70%"

# Row 4
Set-TextValue "D4" "5.305"
Set-TextValue "E4" "3.54%"

# Row 5
Set-TextValue "D5" "0.07455"
Set-TextValue "E5" "11.54%"

# Row 6
Set-TextValue "D6" "7.736"
Set-TextValue "E6" "5.32%"

# Row 7
Set-TextValue "D7" "3.709"
Set-TextValue "E7" "8.91%"

# Row 8
Set-TextValue "D8" "1.553"
Set-TextValue "E8" "14.67%"

# Row 9
Set-TextValue "D9" "0.9230"
Set-TextValue "E9" "0.56%"

# Row 10
Set-TextValue "D10" "0.01637"
Set-TextValue "E10" "2,435.56%"

# Row 11
Set-TextValue "D11" "0.1668"
Set-TextValue "E11" "5.23%"

# Row 12
Set-TextValue "D12" "0.07509"
Set-TextValue "E12" "13.10%"

# Row 13
Set-TextValue "D13" "0.07985"
Set-TextValue "E13" "4.07%"

# Row 14
Set-TextValue "D14" "0.03077"
Set-TextValue "E14" "3.81%"

# Row 15
Set-TextValue "D15" "0.09859"
Set-TextValue "E15" "9.61%"

# Row 16
Set-TextValue "D16" "0.001533"
Set-TextValue "E16" "-3.74%"

# Row 17
Set-TextValue "D17" "0.04545"
Set-TextValue "E17" "0.95%"

# Row 18
Set-TextValue "D18" "0.006459"
Set-TextValue "E18" "2.68%"

# Row 19
Set-TextValue "D19" "3.475"
Set-TextValue "E19" "0.70%"

# Row 20
Set-TextValue "D20" "2.241"
Set-TextValue "E20" "0.97%"

# Row 21
Set-TextValue "D21" "0.3275"
Set-TextValue "E21" "1.92%"

# Row 22
Set-TextValue "D22" "0.1327"
Set-TextValue "E22" "1.35%"

# Row 23
Set-TextValue "D23" "4.218"
Set-TextValue "E23" "3.78%"

# Row 24
Set-TextValue "D24" "0.1620"
Set-TextValue "E24" "4.58%"

# Row 25
Set-TextValue "E25" "0.61%"

# Row 26
Set-TextValue "D26" "0.004530"
Set-TextValue "E26" "9.60%"

# Row 27
Set-TextValue "D27" "0.0001169"
Set-TextValue "E27" "-6.33%"

# Row 28
Set-TextValue "D28" "0.0001664"
Set-TextValue "E28" "2.97%"

# Row 40
Set-TextValue "D40" "0.04497"
Set-TextValue "E40" "6.50%"

# Row 41
Set-TextValue "D41" "0.007307"
Set-TextValue "E41" "8.27%"

# Row 42
Set-TextValue "D42" "0.1366"
Set-TextValue "E42" "9.83%"

# Row 43
Set-TextValue "D43" "0.002259"
Set-TextValue "E43" "14.22%"

# Row 44
Set-TextValue "D44" "0.01384"
Set-TextValue "E44" "18.23%"

# Row 45
Set-TextValue "D45" "0.00006035"
Set-TextValue "E45" "7.52%"

# Row 46
Set-TextValue "E46" "-4.12%"

# Row 47
Set-TextValue "D47" "0.01299"
Set-TextValue "E47" "-0.47%"

